$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 808, shifting existing rows 808:835 down to 809:836
$ws.Rows.Item(808).Insert()

# Populate the new row 808 with the new data record
$ws.Cells.Item(808, 1).Value2 = 6
$ws.Cells.Item(808, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(808, 3).Value2 = "Metropolitana"
$ws.Cells.Item(808, 4).Value2 = 44509
$ws.Cells.Item(808, 5).Value2 = 13
$ws.Cells.Item(808, 6).Value2 = "Fruta"
$ws.Cells.Item(808, 7).Value2 = 100103
$ws.Cells.Item(808, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(808, 9).Value2 = 100103006
$ws.Cells.Item(808, 10).Value2 = "Nectarín"
$ws.Cells.Item(808, 11).Value2 = "Early Glo"
$ws.Cells.Item(808, 12).Value2 = "Segunda"
$ws.Cells.Item(808, 13).Value2 = 12
$ws.Cells.Item(808, 14).Value2 = 550000
$ws.Cells.Item(808, 15).Value2 = 550000
$ws.Cells.Item(808, 16).Value2 = 550000
$ws.Cells.Item(808, 17).Value2 = "`$/bins (420 kilos)"
$ws.Cells.Item(808, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(808, 19).Value2 = 1310
$ws.Cells.Item(808, 20).Value2 = 420
